$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 9 (MIELE Compact C2 Classic), shifting subsequent rows up.
$ws.Rows("9").Delete()

# Update the timestamp column (O) for every remaining data row (2 through 35)
# to reflect the new crawl time.
$ws.Range("O2:O35").Value = "2022-09-11 20:59:46"

Write-Host "Done"
